# Populate the TCR fastq template with example/test data.
# (commit: "add some tcr test data")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TCR")

# ---------------------------------------------------------------------
# 1. Fill in the example values for the "metadata" block (rows 2-7,
#    column C) next to the existing labels in column B.
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "test_prism_trial_id"
$ws.Range("C3").Value = "DFCI"
$ws.Range("C4").Value = "Transcriptome capture"
$ws.Range("C5").Value = "Illumina - TruSeq Stranded PolyA mRNA"
$ws.Range("C6").Value = "Illumina - HiSeq 3000"
$ws.Range("C7").Value = "Paired"

# Give the new value cells (C2:C7) the same "header-ish" look already
# used for the B2:B7 labels: bold black text, light-blue fill, right
# aligned, wrapped, indented, boxed.
$fmt = $ws.Range("C2:C7")
$fmt.HorizontalAlignment = -4152   # xlRight
$fmt.VerticalAlignment = -4108     # xlCenter
$fmt.WrapText = $true
$fmt.IndentLevel = 1
$fmt.Font.Bold = $true
$fmt.Font.Color = 0
$fmt.Interior.Pattern = 1          # xlSolid
$fmt.Interior.Color = 16175794     # RGB(178,210,246) == #B2D2F6

# ---------------------------------------------------------------------
# 2. Add two sample rows of TCR fastq sample data.
# ---------------------------------------------------------------------
$ws.Range("B11").Value = "CTTTPP700.00"
$ws.Range("C11").Value = "/local/path/to/fwd.1.1.1.fastq.gz,/local/path/to/fwd.1.1.1_2.fastq.gz"
$ws.Range("D11").Value = "/local/path/to/rev.1.1.1.fastq.gz"
$ws.Range("E11").Value = 600
$ws.Range("F11").Value = 0.7
$ws.Range("G11").Value = 8
$ws.Range("I11").Value = 1

$ws.Range("B12").Value = "CTTTPP701.00"
$ws.Range("C12").Value = "/local/path/to/fwd.1.2.1.fastq.gz,/local/path/to/fwd.1.2.1_2.fastq.gz"
$ws.Range("D12").Value = "/local/path/to/rev.1.2.1.fastq.gz"
$ws.Range("E12").Value = 650
$ws.Range("F12").Value = 0.8
$ws.Range("G12").Value = 9
$ws.Range("H12").Value = 9
$ws.Range("I12").Value = 1

# ---------------------------------------------------------------------
# 3. Replace the dropdown validations that referenced the "Data
#    Dictionary" sheet ranges with inline literal lists.
# ---------------------------------------------------------------------
foreach ($ref in @("C3","C4","C5","C6","C7")) {
    $ws.Range($ref).Validation.Delete()
}

$ws.Range("C7").Validation.Add(3, 1, 1, '"Paired,Single"')
$ws.Range("C6").Validation.Add(3, 1, 1, '"Illumina - HiSeq 2500,Illumina - HiSeq 3000,Illumina - NextSeq 550,Illumina - HiSeq 4000,Illumina - NovaSeq 6000"')
$ws.Range("C5").Validation.Add(3, 1, 1, '"Agilent,Twist,IDT,NEB,Illumina - TruSeq Stranded PolyA mRNA"')
$ws.Range("C4").Validation.Add(3, 1, 1, '"PolyA capture,Transcriptome capture,Ribo minus"')
$ws.Range("C3").Validation.Add(3, 1, 1, '"DFCI,Mount Sinai,Stanford,MD Anderson"')

# ---------------------------------------------------------------------
# 4. Leave the cursor parked where the author last left it.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("F12").Select()
